$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = 12152.4445045397
$ws.Range("C2").Value = 10824.1750826474
$ws.Range("D2").Value = 15723.86
$ws.Range("F2").Value = 125.732586201085
$ws.Range("B3").Value = 4687.41585855827
$ws.Range("C3").Value = 7333.4968906617
$ws.Range("F3").Value = 126.630048707478
$ws.Range("B4").Value = 4605.52661147761
$ws.Range("C4").Value = 7303.43379856893
$ws.Range("F4").Value = 125.319832057573
$ws.Range("B5").Value = 11869.3607383519
$ws.Range("C5").Value = 11295.4029819199
$ws.Range("F5").Value = 300.822506983837
$ws.Range("B6").Value = 12227.5080442205
$ws.Range("C6").Value = 11909.8138979145
$ws.Range("F6").Value = 336.421411620611
$ws.Range("B7").Value = 12755.9525582825
$ws.Range("C7").Value = 12188.782975305
$ws.Range("F7").Value = 365.438171353057
$ws.Range("B8").Value = 12755.9525582825
$ws.Range("C8").Value = 12186.0555980962
$ws.Range("F8").Value = 365.324530636025
$ws.Range("B9").Value = 12755.9525582825
$ws.Range("C9").Value = 11605.6421522179
$ws.Range("F9").Value = 341.140637057763
$ws.Range("B10").Value = 5130.41124360036
$ws.Range("C10").Value = 8576.28086476087
$ws.Range("F10").Value = 198.53616844102
$ws.Range("B11").Value = 5022.88737964978
$ws.Range("C11").Value = 8751.71967672693
$ws.Range("F11").Value = 205.523429331204
$ws.Range("B12").Value = 12338.8363401026
$ws.Range("C12").Value = 12503.1723355264
$ws.Range("F12").Value = 364.233220170689
$ws.Range("B13").Value = 12338.8363401026
$ws.Range("C13").Value = 12381.1938153886
$ws.Range("F13").Value = 359.150781831615
$ws.Range("B14").Value = 12338.8363401026
$ws.Range("C14").Value = 11942.9338620858
$ws.Range("F14").Value = 340.889950443999
$ws.Range("B15").Value = 12338.8363401026
$ws.Range("C15").Value = 12751.1244041107
$ws.Range("F15").Value = 374.564556361702
